$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the flashcard button link cell to also reference the new image asset
$ws.Range("F8").Value = "/xlsx/Forest.xlsx|/images/flashcard/image.jpeg"

# Move the active selection (matches the recorded view state in the saved file)
$ws.Range("E12").Select()
